$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 106.8
$ws.Range("I39").Value = 106.8
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 320.4
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -24.39999999999998
$ws.Range("N39").ClearContents()
$ws.Range("H62").Value = 4416.1665
$ws.Range("I62").Value = 4324.25
$ws.Range("K62").Value = 4324.25
$ws.Range("M62").Value = -3700.25
$ws.Range("H65").Value = 4416.1665
$ws.Range("I65").Value = 4324.25
$ws.Range("K65").Value = 21621.25
$ws.Range("M65").Value = -18501.25
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H113").Value = 2996.6667
$ws.Range("I113").Value = 2996.6667
$ws.Range("K113").Value = 2996.6667
$ws.Range("M113").Value = 257.3332999999998
$ws.Range("H138").Value = 9094681
$ws.Range("J138").Value = 3498.7778
$ws.Range("L138").Value = 10496.3334
$ws.Range("N138").Value = -20776.3334
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1787.4
$ws.Range("I105").Value = 1313.1666
$ws.Range("J105").Value = 2498.75
$ws.Range("K105").Value = 1313.1666
$ws.Range("L105").Value = 2498.75
$ws.Range("M105").Value = 433.8334
$ws.Range("N105").Value = -5992.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 533.5
$ws.Range("I22").Value = 300.25
$ws.Range("K22").Value = 300.25
$ws.Range("M22").Value = 49.75
$ws.Range("H37").Value = 51
$ws.Range("I37").Value = 51
$ws.Range("K37").Value = 51
$ws.Range("M37").Value = 56
$ws.Range("H58").Value = 10101.8
$ws.Range("I58").Value = 7503.6665
$ws.Range("K58").Value = 7503.6665
$ws.Range("M58").Value = -7300.6665
$ws.Range("H62").Value = 2187
$ws.Range("I62").Value = 2499
$ws.Range("J62").Value = 1875
$ws.Range("K62").Value = 2499
$ws.Range("L62").Value = 1875
$ws.Range("M62").Value = -1875
$ws.Range("N62").Value = -3123
$ws.Range("H65").Value = 2187
$ws.Range("I65").Value = 2499
$ws.Range("J65").Value = 1875
$ws.Range("K65").Value = 12495
$ws.Range("L65").Value = 9375
$ws.Range("M65").Value = -9375
$ws.Range("N65").Value = -15615
$ws.Range("H99").Value = 692.3333
$ws.Range("I99").Value = 688.5
$ws.Range("J99").Value = 700
$ws.Range("K99").Value = 688.5
$ws.Range("L99").Value = 700
$ws.Range("M99").Value = 809.5
$ws.Range("N99").Value = -3696
$ws.Range("H106").Value = 22732.334
$ws.Range("J106").Value = 22732.334
$ws.Range("L106").Value = 22732.334
$ws.Range("N106").Value = -25256.334
$ws.Range("H126").Value = 692.3333
$ws.Range("I126").Value = 688.5
$ws.Range("J126").Value = 700
$ws.Range("K126").Value = 2065.5
$ws.Range("L126").Value = 2100
$ws.Range("M126").Value = 404.5
$ws.Range("N126").Value = -7040
$ws.Range("H136").Value = 10101.8
$ws.Range("I136").Value = 7503.6665
$ws.Range("K136").Value = 22510.9995
$ws.Range("M136").Value = -19960.9995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 191.5
$ws.Range("J40").Value = 162
$ws.Range("L40").Value = 648
$ws.Range("N40").Value = -786
$ws.Range("H104").Value = 4140.3335
$ws.Range("I104").Value = 422
$ws.Range("J104").Value = 5999.5
$ws.Range("K104").Value = 1266
$ws.Range("L104").Value = 17998.5
$ws.Range("M104").Value = 1355
$ws.Range("N104").Value = -23240.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 15025
$ws.Range("J38").Value = 15025
$ws.Range("L38").Value = 15025
$ws.Range("N38").Value = -15951
$ws.Range("H80").Value = 3999.5
$ws.Range("I80").Value = 3999
$ws.Range("K80").Value = 3999
$ws.Range("M80").Value = -3001
$ws.Range("H83").Value = 3999.5
$ws.Range("I83").Value = 3999
$ws.Range("K83").Value = 19995
$ws.Range("M83").Value = -15003
$ws.Range("H102").Value = 3469.5
$ws.Range("I102").Value = 2009.6666
$ws.Range("K102").Value = 2009.6666
$ws.Range("M102").Value = -387.6666
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 3155.111
$ws.Range("I132").Value = 3155.111
$ws.Range("K132").Value = 9465.332999999999
$ws.Range("M132").Value = -6935.332999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1712.1428
$ws.Range("I40").Value = 1527
$ws.Range("J40").Value = 1959
$ws.Range("K40").Value = 1527
$ws.Range("L40").Value = 1959
$ws.Range("M40").Value = -1391
$ws.Range("N40").Value = -2231
$ws.Range("H68").Value = 3192.5
$ws.Range("I68").Value = 3192.5
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3192.5
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2443.5
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 3192.5
$ws.Range("I71").Value = 3192.5
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 15962.5
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -12218.5
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 3967.5
$ws.Range("I122").Value = 3967.5
$ws.Range("K122").Value = 11902.5
$ws.Range("M122").Value = -9452.5
$ws.Range("H132").Value = 7347.5557
$ws.Range("I132").Value = 7347.5557
$ws.Range("K132").Value = 22042.6671
$ws.Range("M132").Value = -19512.6671
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 646.5714
$ws.Range("I122").Value = 646.5714
$ws.Range("K122").Value = 1939.7142
$ws.Range("M122").Value = 510.2857999999999
